$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.660188000000001
$ws.Range("H2").Value = 19.980564
$ws.Range("I2").Value = 0.1500148400131262
$ws.Range("J2").Value = 0.1500148400131261
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 74.08958833333334
$ws.Range("N2").Value = 222.268765
$ws.Range("O2").Value = 0.749631794917355
$ws.Range("P2").Value = 0.7496317949173549
$ws.Range("Q2").Value = 493.4505871426068
$ws.Range("R2").Value = 4441.05528428346
$ws.Range("S2").Value = 0.1124558937832796
$ws.Range("T2").Value = 0.1124558937832796

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.660188000000001
$ws.Range("H3").Value = 19.980564
$ws.Range("I3").Value = 0.1500148400131262
$ws.Range("J3").Value = 0.1500148400131261
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.2521723333333333
$ws.Range("N3").Value = 0.756517
$ws.Range("O3").Value = 0.002551457001146754
$ws.Range("P3").Value = 0.002551457001146754
$ws.Range("Q3").Value = 1.679515148398667
$ws.Range("R3").Value = 15.115636335588
$ws.Range("S3").Value = 0.000382756413827401
$ws.Range("T3").Value = 0.0003827564138274009

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.660188000000001
$ws.Range("H4").Value = 19.980564
$ws.Range("I4").Value = 0.1500148400131262
$ws.Range("J4").Value = 0.1500148400131261
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 21.48095
$ws.Range("N4").Value = 64.44284999999999
$ws.Range("O4").Value = 0.2173423211987967
$ws.Range("P4").Value = 0.2173423211987967
$ws.Range("Q4").Value = 143.0671654186
$ws.Range("R4").Value = 1287.6044887674
$ws.Range("S4").Value = 0.03260457354271896
$ws.Range("T4").Value = 0.03260457354271896

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.660188000000001
$ws.Range("H5").Value = 19.980564
$ws.Range("I5").Value = 0.1500148400131262
$ws.Range("J5").Value = 0.1500148400131261
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.011929
$ws.Range("N5").Value = 9.035787000000001
$ws.Range("O5").Value = 0.03047442688270168
$ws.Range("P5").Value = 0.03047442688270168
$ws.Range("Q5").Value = 20.06001338265201
$ws.Range("R5").Value = 180.540120443868
$ws.Range("S5").Value = 0.004571616273300204
$ws.Range("T5").Value = 0.004571616273300202

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.410331
$ws.Range("H6").Value = 49.230993
$ws.Range("I6").Value = 0.3696281815959916
$ws.Range("J6").Value = 0.3696281815959916
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 74.08958833333334
$ws.Range("N6").Value = 222.268765
$ws.Range("O6").Value = 0.749631794917355
$ws.Range("P6").Value = 0.7496317949173549
$ws.Range("Q6").Value = 1215.834668203738
$ws.Range("R6").Value = 10942.51201383364
$ws.Range("S6").Value = 0.2770850372218412
$ws.Range("T6").Value = 0.2770850372218412

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.410331
$ws.Range("H7").Value = 49.230993
$ws.Range("I7").Value = 0.3696281815959916
$ws.Range("J7").Value = 0.3696281815959916
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.2521723333333333
$ws.Range("N7").Value = 0.756517
$ws.Range("O7").Value = 0.002551457001146754
$ws.Range("P7").Value = 0.002551457001146754
$ws.Range("Q7").Value = 4.138231459042333
$ws.Range("R7").Value = 37.244083131381
$ws.Range("S7").Value = 0.0009430904117542367
$ws.Range("T7").Value = 0.0009430904117542366

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 16.410331
$ws.Range("H8").Value = 49.230993
$ws.Range("I8").Value = 0.3696281815959916
$ws.Range("J8").Value = 0.3696281815959916
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.48095
$ws.Range("N8").Value = 64.44284999999999
$ws.Range("O8").Value = 0.2173423211987967
$ws.Range("P8").Value = 0.2173423211987967
$ws.Range("Q8").Value = 352.5094996944499
$ws.Range("R8").Value = 3172.58549725005
$ws.Range("S8").Value = 0.08033584696856316
$ws.Range("T8").Value = 0.08033584696856316

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 16.410331
$ws.Range("H9").Value = 49.230993
$ws.Range("I9").Value = 0.3696281815959916
$ws.Range("J9").Value = 0.3696281815959916
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.011929
$ws.Range("N9").Value = 9.035787000000001
$ws.Range("O9").Value = 0.03047442688270168
$ws.Range("P9").Value = 0.03047442688270168
$ws.Range("Q9").Value = 49.426751838499
$ws.Range("R9").Value = 444.840766546491
$ws.Range("S9").Value = 0.01126420699383303
$ws.Range("T9").Value = 0.01126420699383303

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.478895
$ws.Range("H10").Value = 43.436685
$ws.Range("I10").Value = 0.3261242951387937
$ws.Range("J10").Value = 0.3261242951387937
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 74.08958833333334
$ws.Range("N10").Value = 222.268765
$ws.Range("O10").Value = 0.749631794917355
$ws.Range("P10").Value = 0.7496317949173549
$ws.Range("Q10").Value = 1072.735370071558
$ws.Range("R10").Value = 9654.618330644025
$ws.Range("S10").Value = 0.2444731407310511
$ws.Range("T10").Value = 0.2444731407310511

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 14.478895
$ws.Range("H11").Value = 43.436685
$ws.Range("I11").Value = 0.3261242951387937
$ws.Range("J11").Value = 0.3261242951387937
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.2521723333333333
$ws.Range("N11").Value = 0.756517
$ws.Range("O11").Value = 0.002551457001146754
$ws.Range("P11").Value = 0.002551457001146754
$ws.Range("Q11").Value = 3.651176736238333
$ws.Range("R11").Value = 32.860590626145
$ws.Range("S11").Value = 0.0008320921160759256
$ws.Range("T11").Value = 0.0008320921160759253

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 14.478895
$ws.Range("H12").Value = 43.436685
$ws.Range("I12").Value = 0.3261242951387937
$ws.Range("J12").Value = 0.3261242951387937
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 21.48095
$ws.Range("N12").Value = 64.44284999999999
$ws.Range("O12").Value = 0.2173423211987967
$ws.Range("P12").Value = 0.2173423211987967
$ws.Range("Q12").Value = 311.0204195502499
$ws.Range("R12").Value = 2799.18377595225
$ws.Range("S12").Value = 0.07088061130478687
$ws.Range("T12").Value = 0.07088061130478686

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 14.478895
$ws.Range("H13").Value = 43.436685
$ws.Range("I13").Value = 0.3261242951387937
$ws.Range("J13").Value = 0.3261242951387937
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.011929
$ws.Range("N13").Value = 9.035787000000001
$ws.Range("O13").Value = 0.03047442688270168
$ws.Range("P13").Value = 0.03047442688270168
$ws.Range("Q13").Value = 43.609403738455
$ws.Range("R13").Value = 392.484633646095
$ws.Range("S13").Value = 0.009938450986879793
$ws.Range("T13").Value = 0.00993845098687979

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.847447
$ws.Range("H14").Value = 20.542341
$ws.Range("I14").Value = 0.1542326832520885
$ws.Range("J14").Value = 0.1542326832520885
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 74.08958833333334
$ws.Range("N14").Value = 222.268765
$ws.Range("O14").Value = 0.749631794917355
$ws.Range("P14").Value = 0.7496317949173549
$ws.Range("Q14").Value = 507.3245293643184
$ws.Range("R14").Value = 4565.920764278865
$ws.Range("S14").Value = 0.1156177231811829
$ws.Range("T14").Value = 0.1156177231811829

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.847447
$ws.Range("H15").Value = 20.542341
$ws.Range("I15").Value = 0.1542326832520885
$ws.Range("J15").Value = 0.1542326832520885
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.2521723333333333
$ws.Range("N15").Value = 0.756517
$ws.Range("O15").Value = 0.002551457001146754
$ws.Range("P15").Value = 0.002551457001146754
$ws.Range("Q15").Value = 1.726736687366333
$ws.Range("R15").Value = 15.540630186297
$ws.Range("S15").Value = 0.0003935180594891908
$ws.Range("T15").Value = 0.0003935180594891908

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.847447
$ws.Range("H16").Value = 20.542341
$ws.Range("I16").Value = 0.1542326832520885
$ws.Range("J16").Value = 0.1542326832520885
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 21.48095
$ws.Range("N16").Value = 64.44284999999999
$ws.Range("O16").Value = 0.2173423211987967
$ws.Range("P16").Value = 0.2173423211987967
$ws.Range("Q16").Value = 147.08966663465
$ws.Range("R16").Value = 1323.80699971185
$ws.Range("S16").Value = 0.03352128938272768
$ws.Range("T16").Value = 0.03352128938272768

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.847447
$ws.Range("H17").Value = 20.542341
$ws.Range("I17").Value = 0.1542326832520885
$ws.Range("J17").Value = 0.1542326832520885
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.011929
$ws.Range("N17").Value = 9.035787000000001
$ws.Range("O17").Value = 0.03047442688270168
$ws.Range("P17").Value = 0.03047442688270168
$ws.Range("Q17").Value = 20.624024195263
$ws.Range("R17").Value = 185.616217757367
$ws.Range("S17").Value = 0.004700152628688658
$ws.Range("T17").Value = 0.004700152628688658
